$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D20").Value = "[파이썬 간단한 게임 만들기] 10. 테트리스(Tetris) - 3탄. 구현"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/559"

$ws.Range("D21").Value = "화자 분리(2) : x-vector 클러스터링 기반의 Speaker Diarization"
$ws.Range("E21").Value = "https://ms-review.tistory.com/11"

$ws.Range("D28").Value = "임피던스 제어 :: ""Simple"" Impedance Control"
$ws.Range("E28").Value = "https://ropiens.tistory.com/129"

$ws.Range("D29").Value = "[만화] 인턴일기 28~33"

$ws.Range("D51").Value = "[javascript] 몫과 나머지 구하기"
$ws.Range("E51").Value = "https://bskyvision.com/816"
